# Infographics.pptx slide 9 updates: octave-up/down labels, grain envelope
# attack wording, and a merged "Right: Fast attack-release linear" run.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)

# --- Shape "TextBox 33" (id 34): Live/Sample/Dry -> Oct Up / Oct Dn / Up/Dn ---
$shp34 = $s.Shapes.Item("TextBox 33")

# Nudge + resize the (rotated) box to its new footprint.
$shp34.Left = 4879954 / 12700
$shp34.Width = 529312 / 12700

$tr34 = $shp34.TextFrame.TextRange

$para1 = $tr34.Paragraphs(1)
$para1.Text = "Oct Up"

$para2 = $tr34.Paragraphs(2)
$para2.Text = "Oct "
$para2.InsertAfter("Dn")

$para3 = $tr34.Paragraphs(3)
$para3.Text = "Up/"
$para3.InsertAfter("Dn")

# --- Shape "TextBox 40" (id 41): Wide -> Slow At, AR -> Fast At ---
$shp41 = $s.Shapes.Item("TextBox 40")
$tr41 = $shp41.TextFrame.TextRange
$tr41.Paragraphs(2).Text = "Slow At"
$tr41.Paragraphs(3).Text = "Fast At"

# --- Shape "TextBox 42" (id 43): Sample Mode -> Oct Mode ---
$shp43 = $s.Shapes.Item("TextBox 42")

$shp43.Left = 4814940 / 12700
$shp43.Width = 750526 / 12700

$tr43 = $shp43.TextFrame.TextRange
# Replace the leading "Sample " with "Oct ", leaving the trailing "Mode" run intact.
$tr43.Characters(1, 7).Text = "Oct "

# --- Shape "Rectangle 57" (id 58): merge the split "Right: Fast ..." runs ---
$shp58 = $s.Shapes.Item("Rectangle 57")
$tr58 = $shp58.TextFrame.TextRange
$para4 = $tr58.Paragraphs(4)
# Force a real change first so the three runs collapse into a single run.
$para4.Text = "__tmp__"
$para4.Text = "  Right: Fast attack-release linear"
